$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Client name updates in column A ---------------------------------------
# "Universal Serv" (row 14) is being replaced across rows 13-14:
#   row 13: AlliedUniversal -> ClearwaterPaper
#   row 14: Universal Serv  -> AlliedUniversal
$ws.Range("A13").Value = "ClearwaterPaper"
$ws.Range("A14").Value = "AlliedUniversal"

# --- C13: boolean FALSE -> literal text "FALSE" -----------------------------
# Assigning the string "FALSE" directly gets auto-coerced back into a Boolean
# by Excel, so build the text value with a formula and paste only the value
# (not the formula) into C13. This preserves C13's existing style (s="3")
# while switching the stored cell type from boolean to shared string.
$ws.Range("Z1").Formula = "=TEXT(FALSE,""@"")"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1").Clear() | Out-Null

# --- Client number update ----------------------------------------------------
$ws.Range("B37").Value = 7594

# --- View / selection update -------------------------------------------------
$ws.Range("B38").Select()
$excel.ActiveWindow.ScrollRow = 8
